$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Remove the two temporary duplicate rows that were added for a one-off
# adjustment ("Sognogfjordane" at row 15, "Moreogromsdal" at row 17).
# Deleting the higher row first keeps the lower row index valid.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(15).Delete()
